$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46-66 down to 47-67
$ws.Rows.Item(46).Insert(-4121)

# Populate the newly inserted row 46 with the new record's data
$ws.Cells.Item(46, 1).Value = 7
$ws.Cells.Item(46, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(46, 3).Value = "Ñuble"
$ws.Cells.Item(46, 4).Value = 44781
$ws.Cells.Item(46, 5).Value = 16
$ws.Cells.Item(46, 6).Value = 100112040
$ws.Cells.Item(46, 7).Value = "Cilantro"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 200
$ws.Cells.Item(46, 11).Value = 700
$ws.Cells.Item(46, 12).Value = 800
$ws.Cells.Item(46, 13).Value = 750
$ws.Cells.Item(46, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(46, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(46, 16).Value = 750
$ws.Cells.Item(46, 17).Value = 1
$ws.Cells.Item(46, 18).Value = "Hortaliza"
